$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 36.75793933333333
$ws.Range("H2").Value = 110.273818
$ws.Range("I2").Value = 0.9858943139827973
$ws.Range("J2").Value = 0.9858943139827971
$ws.Range("M2").Value = 31.61061466666667
$ws.Range("N2").Value = 94.831844
$ws.Range("O2").Value = 0.8860472269592234
$ws.Range("P2").Value = 0.8860472269592234
$ws.Range("Q2").Value = 1161.94105620671
$ws.Range("R2").Value = 10457.46950586039
$ws.Range("S2").Value = 0.8735489229793234
$ws.Range("T2").Value = 0.8735489229793234

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 36.75793933333333
$ws.Range("H3").Value = 110.273818
$ws.Range("I3").Value = 0.9858943139827973
$ws.Range("J3").Value = 0.9858943139827971
$ws.Range("O3").Value = 0.04688826274109129
$ws.Range("P3").Value = 0.04688826274109129
$ws.Range("Q3").Value = 61.48814180035645
$ws.Range("R3").Value = 553.3932762032081
$ws.Range("S3").Value = 0.04622687162897336
$ws.Range("T3").Value = 0.04622687162897334

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 36.75793933333333
$ws.Range("H4").Value = 110.273818
$ws.Range("I4").Value = 0.9858943139827973
$ws.Range("J4").Value = 0.9858943139827971
$ws.Range("M4").Value = 2.392593
$ws.Range("N4").Value = 7.177778999999999
$ws.Range("O4").Value = 0.06706451029968528
$ws.Range("P4").Value = 0.06706451029968527
$ws.Range("Q4").Value = 87.94678834335799
$ws.Range("R4").Value = 791.5210950902219
$ws.Range("S4").Value = 0.06611851937450046
$ws.Range("T4").Value = 0.06611851937450045

# Row 5
$ws.Range("I5").Value = 0.001251989679428792
$ws.Range("J5").Value = 0.001251989679428792
$ws.Range("M5").Value = 31.61061466666667
$ws.Range("N5").Value = 94.831844
$ws.Range("O5").Value = 0.8860472269592234
$ws.Range("P5").Value = 0.8860472269592234
$ws.Range("Q5").Value = 1.475551882025333
$ws.Range("R5").Value = 13.279966938228
$ws.Range("S5").Value = 0.001109321983639448
$ws.Range("T5").Value = 0.001109321983639448

# Row 6
$ws.Range("I6").Value = 0.001251989679428792
$ws.Range("J6").Value = 0.001251989679428792
$ws.Range("O6").Value = 0.04688826274109129
$ws.Range("P6").Value = 0.04688826274109129
$ws.Range("S6").Value = 0.00005870362103819187
$ws.Range("T6").Value = 0.00005870362103819184

# Row 7
$ws.Range("I7").Value = 0.001251989679428792
$ws.Range("J7").Value = 0.001251989679428792
$ws.Range("M7").Value = 2.392593
$ws.Range("N7").Value = 7.177778999999999
$ws.Range("O7").Value = 0.06706451029968528
$ws.Range("P7").Value = 0.06706451029968527
$ws.Range("Q7").Value = 0.111683848647
$ws.Range("R7").Value = 1.005154637823
$ws.Range("S7").Value = 0.00008396407475115192
$ws.Range("T7").Value = 0.00008396407475115187

# Row 8
$ws.Range("G8").Value = 0.4792353333333333
$ws.Range("H8").Value = 1.437706
$ws.Range("I8").Value = 0.01285369633777395
$ws.Range("J8").Value = 0.01285369633777395
$ws.Range("M8").Value = 31.61061466666667
$ws.Range("N8").Value = 94.831844
$ws.Range("O8").Value = 0.8860472269592234
$ws.Range("P8").Value = 0.8860472269592234
$ws.Range("Q8").Value = 15.14892345665155
$ws.Range("R8").Value = 136.340311109864
$ws.Range("S8").Value = 0.01138898199626054
$ws.Range("T8").Value = 0.01138898199626054

# Row 9
$ws.Range("G9").Value = 0.4792353333333333
$ws.Range("H9").Value = 1.437706
$ws.Range("I9").Value = 0.01285369633777395
$ws.Range("J9").Value = 0.01285369633777395
$ws.Range("O9").Value = 0.04688826274109129
$ws.Range("P9").Value = 0.04688826274109129
$ws.Range("Q9").Value = 0.8016578368151112
$ws.Range("R9").Value = 7.214920531336
$ws.Range("S9").Value = 0.000602687491079748
$ws.Range("T9").Value = 0.0006026874910797478

# Row 10
$ws.Range("G10").Value = 0.4792353333333333
$ws.Range("H10").Value = 1.437706
$ws.Range("I10").Value = 0.01285369633777395
$ws.Range("J10").Value = 0.01285369633777395
$ws.Range("M10").Value = 2.392593
$ws.Range("N10").Value = 7.177778999999999
$ws.Range("O10").Value = 0.06706451029968528
$ws.Range("P10").Value = 0.06706451029968527
$ws.Range("Q10").Value = 1.146615103886
$ws.Range("R10").Value = 10.319535934974
$ws.Range("S10").Value = 0.0008620268504336683
$ws.Range("T10").Value = 0.0008620268504336679
